$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking text (e.g. "586.39")
# must be forced to Text format first, otherwise Excel auto-converts them
# to actual numbers (losing the original text formatting / introducing
# floating-point rounding noise).
$textCells = @("D5","D6","D11","D13","D14","D17","D20","D21","D23","D24","D25","D27","D32","D33","D34","D39","D41","D43","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.933.36"
$ws.Range("E2").Value = "  +6.22%  "
$ws.Range("D3").Value = "3.104.88"
$ws.Range("E3").Value = "  +3.64%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "586.39"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.095.96"
$ws.Range("E8").Value = "  +3.82%  "
$ws.Range("E10").Value = "  +10.53%  "
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  +10.75%  "
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  +6.98%  "
$ws.Range("D14").Value = "35.63"
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "3.618.73"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "7.25"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "62.879.88"
$ws.Range("E18").Value = "  +6.13%  "
$ws.Range("D19").Value = "3.107.64"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").Value = "453.21"
$ws.Range("E20").Value = "  +5.34%  "
$ws.Range("D21").Value = "14.10"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "7.59"
$ws.Range("E23").Value = "  +6.71%  "
$ws.Range("D24").Value = "13.60"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "82.03"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "2.25"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("E28").Value = "  +6.56%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  +12.86%  "
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  +10.92%  "
$ws.Range("D33").Value = "26.91"
$ws.Range("E33").Value = "  +4.75%  "
$ws.Range("D34").Value = "2.36"
$ws.Range("E34").Value = "  +12.90%  "
$ws.Range("D35").Value = "0.0₃0807"
$ws.Range("E35").Value = "  +6.89%  "
$ws.Range("E36").Value = "  +4.34%  "
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +13.22%  "
$ws.Range("D39").Value = "51.37"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").Value = "425.29"
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("D42").Value = "2.963.50"
$ws.Range("E42").Value = "  +6.92%  "
$ws.Range("D43").Value = "0.0372"
$ws.Range("E43").Value = "  +5.73%  "
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("E45").Value = "  +9.56%  "
$ws.Range("E46").Value = "  +8.28%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "124.78"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "34.54"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "24.93"
$ws.Range("E51").Value = "  +6.32%  "
